# Refresh cryptos list: Price (col D) and Volume(1h) (col E) for rows 2-51.
# Values are stored as plain text (matching the source feed formatting,
# e.g. "26.106.15", "1.004", "  -1.48%  "), so numeric-looking prices are
# forced to Text before assignment and the temporary number format is
# cleared again afterwards so cell styling is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.106.15'
$ws.Range("E2").Value = '  -1.48%  '
$ws.Range("D3").Value = '1.656.06'
$ws.Range("E3").Value = '  -1.33%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.17'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5157'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.86%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.004'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2626'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.70%  '
$ws.Range("E9").Value = '  -2.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.70'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -5.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07710'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.17%  '
$ws.Range("D12").Value = '1.662.52'
$ws.Range("E12").Value = '  -1.61%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.417'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.13%  '
$ws.Range("D14").Value = '1.883.21'
$ws.Range("E14").Value = '  -1.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5406'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -3.38%  '
$ws.Range("D16").Value = '0.0₅8093'
$ws.Range("E16").Value = '  -3.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.61'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.75%  '
$ws.Range("D18").Value = '26.140.34'
$ws.Range("E18").Value = '  -1.52%  '
$ws.Range("E19").Value = '  +0.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.608'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.99%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '191.09'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.11%  '
$ws.Range("E22").Value = '  -2.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.012'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -5.21%  '
$ws.Range("E24").Value = '  +0.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '139.85'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.78%  '
$ws.Range("E26").Value = '  -4.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.162'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -3.43%  '
$ws.Range("E28").Value = '  -1.34%  '
$ws.Range("E29").Value = '  -2.96%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05966'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -5.01%  '
$ws.Range("E31").Value = '  -1.78%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.538'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.94%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.252'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -4.59%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.600'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -5.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9636'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.93%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.426'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.770'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5660'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -8.46%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01591'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.959'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.74%  '
$ws.Range("E41").Value = '  -0.96%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.003'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.32%  '
$ws.Range("D43").Value = '1.007.87'
$ws.Range("E43").Value = '  -8.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.35'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.28%  '
$ws.Range("D45").Value = '1.797.66'
$ws.Range("E45").Value = '  -1.45%  '
$ws.Range("D46").Value = '0.0₈109'
$ws.Range("E46").Value = '  -2.74%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.57'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.76%  '
$ws.Range("E48").Value = '  +0.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.010'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.97%  '
$ws.Range("E50").Value = '  -0.48%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4196'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.90%  '
